# Weekly update: insert 3 new price rows at the top of the
# "Terminal La Palmera de La Serena - Kiwi" data block (rows 508-510),
# pushing the existing rows 508-519 down to 511-522.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows before the current row 508 (shifts 508:519 -> 511:522)
$ws.Rows.Item(508).Resize(3).Insert()

# New row 508
$ws.Range("A508").Value = 8
$ws.Range("B508").Value = "Terminal La Palmera de La Serena"
$ws.Range("C508").Value = "Coquimbo"
$ws.Range("D508").Value = 45041
$ws.Range("E508").Value = 4
$ws.Range("F508").Value = "Fruta"
$ws.Range("G508").Value = 100101
$ws.Range("H508").Value = "Berries"
$ws.Range("I508").Value = 100101007
$ws.Range("J508").Value = "Kiwi"
$ws.Range("K508").Value = "Hayward"
$ws.Range("L508").Value = "Especial"
$ws.Range("M508").Value = 14
$ws.Range("N508").Value = 400000
$ws.Range("O508").Value = 410000
$ws.Range("P508").Value = 405000
$ws.Range("Q508").Value = "`$/bins (450 kilos)"
$ws.Range("R508").Value = "Región de O'Higgins"
$ws.Range("S508").Value = 900
$ws.Range("T508").Value = 450

# New row 509
$ws.Range("A509").Value = 8
$ws.Range("B509").Value = "Terminal La Palmera de La Serena"
$ws.Range("C509").Value = "Coquimbo"
$ws.Range("D509").Value = 45041
$ws.Range("E509").Value = 4
$ws.Range("F509").Value = "Fruta"
$ws.Range("G509").Value = 100101
$ws.Range("H509").Value = "Berries"
$ws.Range("I509").Value = 100101007
$ws.Range("J509").Value = "Kiwi"
$ws.Range("K509").Value = "Hayward"
$ws.Range("L509").Value = "Primera"
$ws.Range("M509").Value = 10
$ws.Range("N509").Value = 350000
$ws.Range("O509").Value = 360000
$ws.Range("P509").Value = 355000
$ws.Range("Q509").Value = "`$/bins (450 kilos)"
$ws.Range("R509").Value = "Región de O'Higgins"
$ws.Range("S509").Value = 789
$ws.Range("T509").Value = 450

# New row 510
$ws.Range("A510").Value = 8
$ws.Range("B510").Value = "Terminal La Palmera de La Serena"
$ws.Range("C510").Value = "Coquimbo"
$ws.Range("D510").Value = 45041
$ws.Range("E510").Value = 4
$ws.Range("F510").Value = "Fruta"
$ws.Range("G510").Value = 100101
$ws.Range("H510").Value = "Berries"
$ws.Range("I510").Value = 100101007
$ws.Range("J510").Value = "Kiwi"
$ws.Range("K510").Value = "Hayward"
$ws.Range("L510").Value = "Segunda"
$ws.Range("M510").Value = 8
$ws.Range("N510").Value = 300000
$ws.Range("O510").Value = 310000
$ws.Range("P510").Value = 305000
$ws.Range("Q510").Value = "`$/bins (450 kilos)"
$ws.Range("R510").Value = "Región de O'Higgins"
$ws.Range("S510").Value = 678
$ws.Range("T510").Value = 450
